# Modified errors with Cst6 TF group reps 2 & 3 (incorrect fastq files and reads)
# and deleted the failed libraries of Rep 3 that did not have spike in or full run.
#
# The two failed library rows (bioSampleNumber 6 and 7, original rows 7 and 8)
# are removed entirely; all rows below shift up to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two failed-library rows (rows 7 and 8).
$ws.Rows("7:8").Delete()

# Leave the selection where the user's last edit landed.
$ws.Range("D21").Select() | Out-Null
